$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Composicion de Lote" (I) column had stray empty placeholder cells on
# rows 3-15 (no actual content) -- drop them so only the real header (I1) remains.
$ws.Range("I3:I15").ClearContents()

# Newly appended catalogue rows (16-26)
$ws.Range("A16").Value = "6XS18565"
$ws.Range("B16").Value = "SUAVINEX SPRAY HIGIENIZANTE DE MANOS 100ML"
$ws.Range("C16").Value = "ANEXOS"
$ws.Range("D16").Value = "Tiene PT"
$ws.Range("E16").Value = "Tiene ES"
$ws.Range("F16").Value = "Tiene IT"
$ws.Range("H16").Value = "ML"
$ws.Range("J16").Value = "Solo Revisión"
$ws.Range("G16").Formula = "'100"
$ws.Range("G16").ClearFormats()

$ws.Range("A17").Value = "6XS18552"
$ws.Range("B17").Value = "SUAVINEX LOCION 500ML + GEL 500ML PACK"
$ws.Range("C17").Value = "ANEXOS"
$ws.Range("D17").Value = "Tiene PT"
$ws.Range("E17").Value = "Tiene ES"
$ws.Range("F17").Value = "Tiene IT"
$ws.Range("H17").Value = "UND"
$ws.Range("J17").Value = "Solo Revisión"
$ws.Range("G17").Formula = "'2"
$ws.Range("G17").ClearFormats()

$ws.Range("A18").Value = "6XS18549"
$ws.Range("B18").Value = "SUAVINEX BALSAMO BABY PECTORAL AROMATICO 50ML"
$ws.Range("C18").Value = "ANEXOS"
$ws.Range("D18").Value = "Tiene PT"
$ws.Range("E18").Value = "Tiene ES"
$ws.Range("F18").Value = "Tiene IT"
$ws.Range("H18").Value = "ML"
$ws.Range("J18").Value = "Solo Revisión"
$ws.Range("G18").Formula = "'50"
$ws.Range("G18").ClearFormats()

$ws.Range("A19").Value = "6XS18553"
$ws.Range("B19").Value = "SUAVINEX DEO KIDS ROLL-ON 50ML"
$ws.Range("D19").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E19").Value = "Tiene ES"
$ws.Range("F19").Value = "Tiene IT"
$ws.Range("H19").Value = "ML"
$ws.Range("J19").Value = "Revisado y Traducido"
$ws.Range("G19").Formula = "'50"
$ws.Range("G19").ClearFormats()

$ws.Range("A20").Value = "6XS18562"
$ws.Range("B20").Value = "SUAVINEX LIMPIADOR NASAL"
$ws.Range("C20").Value = "ANEXOS"
$ws.Range("D20").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E20").Value = "Tiene ES"
$ws.Range("F20").Value = "Tiene IT"
$ws.Range("H20").Value = "UND"
$ws.Range("J20").Value = "Revisado y Traducido"
$ws.Range("G20").Formula = "'1"
$ws.Range("G20").ClearFormats()

$ws.Range("A21").Value = "6XS18564"
$ws.Range("B21").Value = "SUAVINEX SPRAY NASAL AGUA DE MAR Y ALOE 120ML +3A"
$ws.Range("C21").Value = "ANEXOS"
$ws.Range("D21").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E21").Value = "Tiene ES"
$ws.Range("F21").Value = "No Tiene IT - TRADOTTO"
$ws.Range("H21").Value = "ML"
$ws.Range("J21").Value = "Revisado y Traducido"
$ws.Range("G21").Formula = "'120"
$ws.Range("G21").ClearFormats()

$ws.Range("A22").Value = "6XS18564"
$ws.Range("B22").Value = "SUAVINEX SPRAY NASAL AGUA DE MAR Y ALOE 120ML +3A"
$ws.Range("C22").Value = "ANEXOS"
$ws.Range("D22").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E22").Value = "Tiene ES"
$ws.Range("F22").Value = "No Tiene IT - TRADOTTO"
$ws.Range("H22").Value = "ML"
$ws.Range("J22").Value = "Revisado y Traducido"
$ws.Range("G22").Formula = "'120"
$ws.Range("G22").ClearFormats()

$ws.Range("A23").Value = "6XS18563"
$ws.Range("B23").Value = "SUAVINEX SPRAY NASAL HIPERTONICO 120ML +3M"
$ws.Range("C23").Value = "ANEXOS"
$ws.Range("D23").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E23").Value = "Tiene ES"
$ws.Range("F23").Value = "Tiene IT"
$ws.Range("H23").Value = "ML"
$ws.Range("J23").Value = "Revisado y Traducido"
$ws.Range("G23").Formula = "'120"
$ws.Range("G23").ClearFormats()

$ws.Range("A24").Value = "6XS18555"
$ws.Range("B24").Value = "SUAVINEX MOM BALSAMO PEZON 30ML"
$ws.Range("C24").Value = "ANEXOS"
$ws.Range("D24").Value = "Tiene PT"
$ws.Range("E24").Value = "Tiene ES"
$ws.Range("F24").Value = "No Tiene IT - TRADOTTO"
$ws.Range("H24").Value = "ML"
$ws.Range("J24").Value = "Solo Revisión"
$ws.Range("G24").Formula = "'30"
$ws.Range("G24").ClearFormats()

$ws.Range("A25").Value = "6XS18556"
$ws.Range("B25").Value = "SUAVINEX MOM ACEITE ESTRIAS 100ML"
$ws.Range("C25").Value = "ANEXOS"
$ws.Range("D25").Value = "Tiene PT"
$ws.Range("E25").Value = "Tiene ES"
$ws.Range("F25").Value = "No Tiene IT - TRADOTTO"
$ws.Range("H25").Value = "ML"
$ws.Range("J25").Value = "Solo Revisión"
$ws.Range("G25").Formula = "'100"
$ws.Range("G25").ClearFormats()

$ws.Range("A26").Value = "0TF27094"
$ws.Range("B26").Value = "THE POTIONS CENTELLA ASIATICA WATER ESSENCE 50ML"
$ws.Range("C26").Value = "TRAT.FEMENINO"
$ws.Range("D26").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E26").Value = "Tiene ES"
$ws.Range("F26").Value = "No Tiene IT - TRADOTTO"
$ws.Range("H26").Value = "ML"
$ws.Range("J26").Value = "Revisado y Traducido"
$ws.Range("G26").Formula = "'50"
$ws.Range("G26").ClearFormats()
